$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (e.g. "1.00" -> 1, "590.70" -> 590.7)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.593.70'
$ws.Range("E2").Value = '  +3.39%  '
$ws.Range("D3").Value = '3.504.37'
$ws.Range("E3").Value = '  +1.76%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '590.70'
$ws.Range("E5").Value = '  +2.83%  '
$ws.Range("D6").Value = '170.26'
$ws.Range("E6").Value = '  +3.19%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.500.44'
$ws.Range("E8").Value = '  +1.65%  '
$ws.Range("D9").Value = '0.592'
$ws.Range("E9").Value = '  +6.19%  '
$ws.Range("D10").Value = '7.34'
$ws.Range("E10").Value = '  +0.46%  '
$ws.Range("E11").Value = '  +4.99%  '
$ws.Range("D12").Value = '0.440'
$ws.Range("E12").Value = '  +2.74%  '
$ws.Range("D13").Value = '4.108.75'
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("E14").Value = '  -0.78%  '
$ws.Range("D15").Value = '28.36'
$ws.Range("E15").Value = '  +3.29%  '
$ws.Range("D16").Value = '0.0000179'
$ws.Range("E16").Value = '  +1.65%  '
$ws.Range("D17").Value = '66.615.13'
$ws.Range("E17").Value = '  +3.34%  '
$ws.Range("D18").Value = '3.532.56'
$ws.Range("E18").Value = '  +2.92%  '
$ws.Range("E19").Value = '  +1.98%  '
$ws.Range("D20").Value = '14.17'
$ws.Range("E20").Value = '  +3.28%  '
$ws.Range("D21").Value = '392.49'
$ws.Range("E21").Value = '  +3.35%  '
$ws.Range("D22").Value = '8.00'
$ws.Range("E22").Value = '  +1.23%  '
$ws.Range("D23").Value = '73.07'
$ws.Range("E23").Value = '  +2.19%  '
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("E25").Value = '  +3.06%  '
$ws.Range("E26").Value = '  +4.15%  '
$ws.Range("D27").Value = '10.37'
$ws.Range("E27").Value = '  +7.60%  '
$ws.Range("E28").Value = '  +2.35%  '
$ws.Range("E29").Value = '  +0.33%  '
$ws.Range("D30").Value = '6.33'
$ws.Range("E30").Value = '  +3.21%  '
$ws.Range("D31").Value = '1.48'
$ws.Range("E31").Value = '  +4.01%  '
$ws.Range("E32").Value = '  +2.39%  '
$ws.Range("D33").Value = '23.65'
$ws.Range("E33").Value = '  +2.50%  '
$ws.Range("E34").Value = '  +3.08%  '
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").Value = '1.63'
$ws.Range("E36").Value = '  +7.09%  '
$ws.Range("E37").Value = '  +1.41%  '
$ws.Range("D38").Value = '0.885'
$ws.Range("E38").Value = '  +2.52%  '
$ws.Range("E39").Value = '  +4.16%  '
$ws.Range("B40").Value = 'InjectiveProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D40").Value = '27.96'
$ws.Range("E40").Value = '  +5.61%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '6.84'
$ws.Range("E41").Value = '  +4.77%  '
$ws.Range("E42").Value = '  +5.42%  '
$ws.Range("D43").Value = '0.0747'
$ws.Range("E43").Value = '  +1.82%  '
$ws.Range("D44").Value = '26.57'
$ws.Range("E44").Value = '  +1.83%  '
$ws.Range("D45").Value = '2.795.28'
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("D46").Value = '43.23'
$ws.Range("E46").Value = '  +0.57%  '
$ws.Range("E47").Value = '  +0.84%  '
$ws.Range("E48").Value = '  -0.59%  '
$ws.Range("D49").Value = '351.08'
$ws.Range("E49").Value = '  +3.47%  '
$ws.Range("E50").Value = '  +2.90%  '
$ws.Range("D51").Value = '33.75'
$ws.Range("E51").Value = '  +11.15%  '
